# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
# Both sheets contain the same event list, and each listed row's F value
# was bumped up as the source data was refreshed.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 140
    3  = 1666
    7  = 11732
    8  = 35
    10 = 467
    11 = 386
    13 = 828
    14 = 13413
    15 = 13270
    16 = 36
    18 = 19
    20 = 262
    23 = 147
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
